$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.429.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.11%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.624.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.64%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "

# Row 6
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.495"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.99%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.56%  "

# Row 8
$ws.Range("E8").Value = "  -0.96%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0838"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.51%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.852.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.37%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.628.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.46%  "

# Row 14
$ws.Range("E14").Value = "  +1.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.97%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.431.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.94%  "

# Row 18
$ws.Range("E18").Value = "  +0.35%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.77%  "

# Row 20
$ws.Range("E20").Value = "  +0.54%  "

# Row 21
$ws.Range("E21").Value = "  -1.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "

# Row 23
$ws.Range("E23").Value = "  -1.41%  "

# Row 24
$ws.Range("E24").Value = "  +4.96%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.89%  "

# Row 26
$ws.Range("E26").Value = "  +0.59%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.120"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.97%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.47%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.96%  "

# Row 30
$ws.Range("E30").Value = "  -2.51%  "

# Row 31
$ws.Range("E31").Value = "  -1.64%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.22%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.18%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.216.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.30%  "

# Row 37
$ws.Range("E37").Value = "  +3.77%  "

# Row 38
$ws.Range("E38").Value = "  +0.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.793"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.69%  "

# Row 40
$ws.Range("E40").Value = "  -0.92%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.94%  "

# Row 42
$ws.Range("E42").Value = "  +0.72%  "

# Row 43
$ws.Range("E43").Value = "  -1.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.762.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.28%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.36%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.71%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.19%  "

# Row 48
$ws.Range("E48").Value = "  -1.87%  "

# Row 49
$ws.Range("E49").Value = "  -0.23%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.54%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.408"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
